$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.344.76"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.092.20"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'582.10"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'144.30"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.084.35"
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("D11").Value = "'5.63"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  -2.63%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'37.64"
$ws.Range("E14").Value = "  +6.16%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "3.606.00"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "63.225.61"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "3.087.29"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "'459.84"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "'14.17"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'7.43"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "'80.98"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "'2.12"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'8.88"
$ws.Range("E28").Value = "  +7.27%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("D32").Value = "'6.78"
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").Value = "'26.65"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").Value = "'2.30"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "  +5.18%  "
$ws.Range("D39").Value = "'5.98"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "'50.19"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "'435.06"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").Value = "'8.70"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").Value = "'0.0367"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "2.856.86"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("E46").Value = "  -3.88%  "
$ws.Range("D47").Value = "'35.96"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "'123.97"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "'24.02"
$ws.Range("E51").Value = "  -2.71%  "
